$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "№ комиссии [number], направление" — the trailing "," run (which
# was underlined, grouped together with "[number]") and the following
# " направление" run (not underlined) get merged into a single, non
# underlined run containing ", направление".
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(", направление", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $commaStart = $rng1.Start
    $dirRng = $d.Range($commaStart + 1, $rng1.End)

    $find = $dirRng.Find
    $find.ClearFormatting()
    $find.Text = " направление"
    $find.Replacement.ClearFormatting()
    $find.Replacement.Text = ", направление"
    $find.Forward = $true
    $find.Wrap = 0
    $find.Format = $false
    $find.MatchCase = $false
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.MatchSoundsLike = $false
    $find.MatchAllWordForms = $false
    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null

    # Remove the now-redundant standalone "," run that used to sit in front of
    # " направление" (its text got duplicated onto the merged run above).
    $d.Range($commaStart, $commaStart + 1).Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# Change 2: insert a line-wrap (newline + indent) between "экзамена" and
# "или" in the long "Оценка СПбГУ ..." sentence, replacing the single space
# that used to separate the two words.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("экзамена или", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $spacePos = $rng2.Start + "экзамена".Length
    $spaceRng = $d.Range($spacePos, $spacePos + 1)
    $spaceRng.Text = "`n                                "
}

# ---------------------------------------------------------------------------
# Change 3: "ФИО  члена" (double space) -> "ФИО члена" (single space).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("ФИО  члена", $true, $false, $false, $false, $false, $true, 1, $false, "ФИО члена", 2) | Out-Null
